# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (and before the
#    "总计" summary sheet) containing the Q1-2022 fund-holding detail rows.
# 2. Rebuild the "总计" summary sheet so it keeps its previous three rows but
#    gains a new leading row for "2022-Q1".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell as literal TEXT (no numeric auto-detect,
# no residual number-format styling left behind).
function Set-TextCell($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Step 1: locate the old "总计" sheet; its previous rows ("2021-Q4" and
# "2021-Q1") are re-declared below (together with the new "2022-Q1" row) so
# the replacement sheet can be repopulated in one go.
# ---------------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# Step 2: delete the existing "总计" sheet. This frees up its sheetId so the
# new "2022-Q1" sheet can take it over (matching the author's sheetId order),
# and we'll recreate "总计" afterwards in the right tab position.
# ---------------------------------------------------------------------------
$oldTotal.Delete() | Out-Null

# ---------------------------------------------------------------------------
# Step 3: create the new "2022-Q1" sheet right after "2021-Q4".
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q1_2022 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4Sheet)
$q1_2022.Name = "2022-Q1"

# Copy the header / index-column formatting (style used for B1:H1 and the
# A-column index cells) from the "2021-Q4" sheet, which already carries the
# bordered, bold, centered style used throughout this workbook.
$q4Sheet.Range("B1:H1").Copy() | Out-Null
$q1_2022.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$q4Sheet.Range("A2").Copy() | Out-Null
$q1_2022.Range("A2:A8").PasteSpecial(-4122) | Out-Null

$q1_2022.Application.CutCopyMode = $false

# Header row.
$q1_2022.Range("B1").Value = "基金代码"
$q1_2022.Range("C1").Value = "基金名称"
$q1_2022.Range("D1").Value = "基金规模"
$q1_2022.Range("E1").Value = "股票总仓位"
$q1_2022.Range("F1").Value = "仓位占比"
$q1_2022.Range("G1").Value = "持有市值(亿元)"
$q1_2022.Range("H1").Value = "仓位排名"

# Detail rows: 基金代码 / 基金名称 / 基金规模 / 股票总仓位 / 仓位占比 /
# 持有市值(亿元) are all stored as literal text; 仓位排名 is numeric; the
# leading index column A is numeric (0-based).
$fundRows = @(
    @("360007", "光大保德信优势配置混合",         "8.68", "85.46", "3.92", "0.3403", 4),
    @("166109", "信达澳银量化先锋混合（LOF）A",     "1.06", "89.44", "2.75", "0.0292", 8),
    @("970046", "东海证券海睿健行灵活配置混合A",   "0.57", "77.45", "2.88", "0.0164", 7),
    @("002810", "金信转型创新成长灵活配置混合",     "0.18", "81.12", "3.62", "0.0065", 6),
    @("970047", "东海证券海睿健行灵活配置混合B",   "0.19", "77.45", "2.88", "0.0055", 7),
    @("166110", "信达澳银量化先锋混合（LOF）C",     "0.09", "89.44", "2.75", "0.0025", 8),
    @("970083", "东海证券海盈6个月持有期混合",     "0.14", "20.43", "1.00", "0.0014", 9)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $q1_2022.Range("A$r").Value = $i

    Set-TextCell $q1_2022.Range("B$r") $row[0]
    Set-TextCell $q1_2022.Range("C$r") $row[1]
    Set-TextCell $q1_2022.Range("D$r") $row[2]
    Set-TextCell $q1_2022.Range("E$r") $row[3]
    Set-TextCell $q1_2022.Range("F$r") $row[4]
    Set-TextCell $q1_2022.Range("G$r") $row[5]

    $q1_2022.Range("H$r").Value = $row[6]
}

# ---------------------------------------------------------------------------
# Step 4: recreate the "总计" sheet after "2022-Q1" and repopulate it with
# the previous rows plus the new leading "2022-Q1" row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1_2022)
$total.Name = "总计"

$q4Sheet.Range("B1:D1").Copy() | Out-Null
$total.Range("B1:D1").PasteSpecial(-4122) | Out-Null

$q4Sheet.Range("A2").Copy() | Out-Null
$total.Range("A2:A4").PasteSpecial(-4122) | Out-Null

$total.Application.CutCopyMode = $false

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

# New "2022-Q1" row first, followed by the sheet's previous rows, unchanged.
$allTotalRows = @(
    @("2022-Q1", 7, 0.4),
    @("2021-Q4", 4, 6.76),
    @("2021-Q1", 2, 0.02)
)

for ($i = 0; $i -lt $allTotalRows.Count; $i++) {
    $r = $i + 2
    $row = $allTotalRows[$i]

    $total.Range("A$r").Value = $i
    $total.Range("B$r").Value = $row[0]
    $total.Range("C$r").Value = $row[1]
    $total.Range("D$r").Value = $row[2]
}

# ---------------------------------------------------------------------------
# Step 5: restore the originally active tab ("2021-Q1"), since adding /
# deleting sheets shifts the active-sheet selection as a side effect.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
